# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Papa" (Asterix, 1a nueva(o)) at
# Vega Monumental Concepción. The new record is placed at row 347, pushing
# the existing rows 347-374 down to 348-375.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 347, shifting rows 347:374 -> 348:375
$ws.Rows.Item(347).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(347, 1).Value()  = 11
$ws.Cells.Item(347, 2).Value()  = "Vega Monumental Concepción"
$ws.Cells.Item(347, 3).Value()  = "Bíobío"
$ws.Cells.Item(347, 4).Value()  = 44931
$ws.Cells.Item(347, 5).Value()  = 8
$ws.Cells.Item(347, 6).Value()  = 100114001
$ws.Cells.Item(347, 7).Value()  = "Papa"
$ws.Cells.Item(347, 8).Value()  = "Asterix"
$ws.Cells.Item(347, 9).Value()  = "1a nueva(o)"
$ws.Cells.Item(347, 10).Value() = 220
$ws.Cells.Item(347, 11).Value() = 12500
$ws.Cells.Item(347, 12).Value() = 13000
$ws.Cells.Item(347, 13).Value() = 12727
$ws.Cells.Item(347, 14).Value() = "$/malla 25 kilos"
$ws.Cells.Item(347, 15).Value() = "Región de La Araucanía"
$ws.Cells.Item(347, 16).Value() = 509
$ws.Cells.Item(347, 17).Value() = 25
$ws.Cells.Item(347, 18).Value() = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D
$ws.Cells.Item(347, 4).NumberFormat() = $ws.Cells.Item(348, 4).NumberFormat()
